$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain text so values serialize
# the same way as the source data (plain strings, not auto-converted
# numbers/percentages).
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Updated Price values (column D)
$ws.Range("D2").Value = "329.65"
$ws.Range("D3").Value = "40.09"
$ws.Range("D4").Value = "5.387"
$ws.Range("D5").Value = "0.08099"
$ws.Range("D6").Value = "4.530"
$ws.Range("D7").Value = "8.654"
$ws.Range("D10").Value = "0.9426"
$ws.Range("D11").Value = "0.1358"
$ws.Range("D12").Value = "0.1978"
$ws.Range("D13").Value = "0.09306"
$ws.Range("D14").Value = "0.03564"
$ws.Range("D15").Value = "0.09584"
$ws.Range("D16").Value = "0.001316"
$ws.Range("D17").Value = "0.006406"
$ws.Range("D18").Value = "3.364"
$ws.Range("D19").Value = "0.3523"
$ws.Range("D20").Value = "7.235"
$ws.Range("D23").Value = "0.04426"
$ws.Range("D24").Value = "0.001220"
$ws.Range("D25").Value = "0.004269"
$ws.Range("D26").Value = "0.0001199"
$ws.Range("D27").Value = "0.0003988"
$ws.Range("D39").Value = "0.02494"
$ws.Range("D40").Value = "0.05228"
$ws.Range("D41").Value = "0.007634"
$ws.Range("D43").Value = "0.009183"
$ws.Range("D44").Value = "0.002169"
$ws.Range("D45").Value = "0.01077"
$ws.Range("D46").Value = "0.00006583"
$ws.Range("D48").Value = "0.002399"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("D51").Value = "0.0001999"

# Updated Volume(1h) values (column E)
$ws.Range("E2").Value = "7.02%"
$ws.Range("E3").Value = "7.78%"
$ws.Range("E4").Value = "5.20%"
$ws.Range("E5").Value = "3.40%"
$ws.Range("E6").Value = "3.05%"
$ws.Range("E7").Value = "4.91%"
$ws.Range("E8").Value = "1.90%"
$ws.Range("E9").Value = "1.16%"
$ws.Range("E10").Value = "2.36%"
$ws.Range("E11").Value = "25.49%"
$ws.Range("E12").Value = "4.35%"
$ws.Range("E13").Value = "5.14%"
$ws.Range("E14").Value = "6.69%"
$ws.Range("E15").Value = "-0.07%"
$ws.Range("E16").Value = "-4.29%"
$ws.Range("E17").Value = "11.66%"
$ws.Range("E18").Value = "-1.41%"
$ws.Range("E19").Value = "2.87%"
$ws.Range("E20").Value = "14.87%"
$ws.Range("E21").Value = "3.51%"
$ws.Range("E22").Value = "6.00%"
$ws.Range("E23").Value = "1.58%"
$ws.Range("E24").Value = "2.25%"
$ws.Range("E25").Value = "-0.04%"
$ws.Range("E26").Value = "-14.33%"
$ws.Range("E27").Value = "-0.10%"
$ws.Range("E39").Value = "14.65%"
$ws.Range("E40").Value = "3.83%"
$ws.Range("E41").Value = "0.81%"
$ws.Range("E42").Value = "5.67%"
$ws.Range("E43").Value = "6.00%"
$ws.Range("E44").Value = "9.59%"
$ws.Range("E45").Value = "36.35%"
$ws.Range("E46").Value = "0.86%"
$ws.Range("E47").Value = "0.02%"
$ws.Range("E48").Value = "139.37%"
$ws.Range("E49").Value = "1.52%"
$ws.Range("E50").Value = "0.02%"
$ws.Range("E51").Value = "0.02%"

# Restore default (unstyled) formatting so these cells keep looking
# like the rest of the unstyled data cells.
$dataRange.Style = "Normal"
